# update arrival data with 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New header row (row 2): YEAR/LOCATION, Date/time, flow/height, NOTES
# ------------------------------------------------------------------
$ws.Range("A2").Value2 = "YEAR/LOCATION"
$ws.Range("B2").Value2 = "Date/time"
$ws.Range("C2").Value2 = "flow/height"
$ws.Range("D2").Value2 = "NOTES"
$ws.Range("A2:D2").Font.Bold = $true

# ------------------------------------------------------------------
# 2) Bold the "year" marker cells that already exist (2021/2022/2023)
# ------------------------------------------------------------------
$ws.Range("A3").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true
$ws.Range("A17").Font.Bold = $true

# ------------------------------------------------------------------
# 3) Fix row 13: B13 used to hold the literal text "6/2/22/9:33";
#    replace it with the real date/time value, formatted like the
#    other Date/time cells in column B.
# ------------------------------------------------------------------
$ws.Range("B13").Value2 = 44714.397916666669
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 4) Append the 2024 arrival data block (rows 24-29)
# ------------------------------------------------------------------
$ws.Range("A24").Value2 = 2024
$ws.Range("A24").Font.Bold = $true

$ws.Range("A25").Value2 = "Caballo"
$ws.Range("B25").Value2 = 45359.354166666664
$ws.Range("C25").Value2 = 505

$ws.Range("B26").Value2 = 45359.395833333336
$ws.Range("C26").Value2 = 2021

$ws.Range("A27").Value2 = "Leasburg"
$ws.Range("B27").Value2 = 45360.271527777775
$ws.Range("C27").Value2 = 341

$ws.Range("B28").Value2 = 45360.334027777775
$ws.Range("C28").Value2 = 1051

$ws.Range("A29").Value2 = "Picacho"
$ws.Range("B29").Value2 = 45360.750694444447
$ws.Range("C29").Value2 = 4.36

# Match the date/time number format used elsewhere in column B for the
# new rows ...
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B25:B28").PasteSpecial(-4122) | Out-Null
# ... except B29, which only shows the date (no time-of-day component).
$ws.Range("B29").NumberFormat = "mm-dd-yy"

# ------------------------------------------------------------------
# 5) Column widths for the now-wider header columns
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(2).ColumnWidth = 21.75
$ws.Columns.Item(3).ColumnWidth = 12.42

# ------------------------------------------------------------------
# 6) Leave selection where the author left it
# ------------------------------------------------------------------
$ws.Range("B30").Select() | Out-Null
